$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "2026年３月期の通期連結業績予想につきましては、米国政府の関税政策の動向などを慎重に見極めるため、未定としておりましたが、現時点で入手可能な情報や予測を踏まえ、上記の通り連結業績予想を算定しましたので公表するものであります。"
